$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prefix the per-category line-item labels with their category name so
# e.g. "     New nominations" under "Civilian " becomes
# "     Civilian, New nominations", matching the other category blocks'
# already-distinct labels.
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Carryover nominations"
$ws.Range("A9").Value  = "     Civilian, Confirmed "
$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("A11").Value = "     Civilian, Returned to White House "

$ws.Range("A13").Value = "     Other Civilian, New nominations"
$ws.Range("A14").Value = "     Other Civilian, Confirmed "
$ws.Range("A15").Value = "     Other Civilian, Returned to White House "

$ws.Range("A17").Value = "     Air Force, New nominations"
$ws.Range("A18").Value = "     Air Force, Confirmed "
$ws.Range("A19").Value = "     Air Force, Withdrawn "
$ws.Range("A20").Value = "     Air Force, Returned to White House "

$ws.Range("A22").Value = "     Army, New nominations"
$ws.Range("A23").Value = "     Army, Confirmed "
$ws.Range("A24").Value = "     Army, Withdrawn "
$ws.Range("A25").Value = "     Army, Returned to White House "

$ws.Range("A27").Value = "     Navy, New nominations"
$ws.Range("A28").Value = "     Navy, Confirmed "
$ws.Range("A29").Value = "     Navy, Returned to White House "

$ws.Range("A31").Value = "     Marine Corps, New nominations"
$ws.Range("A32").Value = "     Marine Corps, Confirmed "
$ws.Range("A33").Value = "     Marine Corps, Returned to White House "

# The Summary section reorders its first two totals: "Total nominations
# received this Session" (old row 36) moves above "Total nominations
# carried over from the First Session" (old row 35). Swap their
# value+style (full Copy, not just .Value) via a scratch cell so each
# number keeps its original number formatting.
$ws.Range("B35").Copy($ws.Range("D1"))
$ws.Range("B36").Copy($ws.Range("B35"))
$ws.Range("D1").Copy($ws.Range("B36"))
$ws.Range("D1").Clear()

# Remove the "Summary" section header row (old row 34); everything below
# shifts up by one.
$ws.Rows("34").Delete()

# Relabel the summary totals (old rows 35-40, now rows 34-39) to their
# new shorter names; values/styles already carried up correctly by the
# swap + row delete above.
$ws.Range("A34").Value = "Total new nominations"
$ws.Range("A35").Value = "Total carryover nominations"
$ws.Range("A36").Value = "Total confirmed "
$ws.Range("A37").Value = "Total unconfirmed "
$ws.Range("A38").Value = "Total withdrawn "
$ws.Range("A39").Value = "Total returned to the White House "
